$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'football youth compression pants'
    'knee guard'
    'knee sleeve basketball pair'
    'soccer tights for boys'
    'knee pads volleyball girls'
    'medias basketball'
    'boys yoga'
    'baseball chart'
    'basketball legs'
    'youth knee pad'
    'softball band'
    'need pads volleyball'
    'mens yoga pants capri'
    'compression shorts knee length'
    'men thigh compression pants'
    'basketball leg compression'
    'boys youth leggings'
    'fabric softball'
    'bump pad'
    'men knee pads for work'
    'black spandex pants men'
    'youth volleyball kneepads'
    'compression tight men'
    'black compression tights men'
    'basketball sleeve with pad'
    'compression padded knee sleeve'
    'compression pants black'
    'men capris'
    'calf pad'
    'knee pads volleyball girls youth'
    'men yoga pants'
    'slim knee pads'
    'hockey leg pads'
    'soccer leggings men'
    'knee pads volleyball large'
    'softball fabric'
    'knee pad thick'
    'yoga for knees'
    'knee pads for boys'
    'knee sleeve honeycomb'
    'black mens leggings'
    'used softballs'
    'youth tights for sports'
    'boy compression leggings'
    'cold compression pants'
    'youth knee guards'
    'paintball knee'
    'knee pad construction'
    'basketball knee sleeve youth'
    'mens wrestling shorts'
    'recovery compression tights men'
    'womens sliding shorts softball'
    'gym kneeling pad'
    'gym leggings for men'
    'durable pants'
    'adult leggings'
    'work wear knee pads'
    'black lacrosse shorts'
    'hex skin padding'
    'sport leggings boys'
    'need pads for construction'
    'football girdle youth'
    '5 inch foam basketball'
    'mens skin tight leggings'
    'silicon knee pads'
    'extra small baseball pants'
    'football pads for men'
    'mens kneepads'
    'thigh protector men'
    'long volleyball knee pads'
    'paintball pants youth'
    'gym leggings men'
    'compression running capris'
    'youth running pants boys'
    'knee sleeve baseball'
    'football pads youth'
    'volleyball spandex shorts'
    'work in baseball'
    'thick yoga knee pad'
    'running tights mens'
    'basketballs 28 5'
    'softball pants'
    'youth volleyball spandex'
    'work pants knee pads'
    'size small baseball pants'
    'boys sport tights'
    'mens medium tall athletic pants'
    'compression sleeve knee youth'
    'compression knee sleeves for basketball'
    'mens knee pads for work'
    'adult knee pads for work'
    'leggings cycling'
    'male athletic tights'
    'knee compression sleeve basketball'
    'compressions knee'
    'youth running tights'
    'knee chart'
    'long softball pants'
    'compression pants men soccer'
    'limber support'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
